$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Administrator, Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat"
$ws.Range("G3").Value = "Administrator, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Veronia Rafat"
$ws.Range("G4").Value = "Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad"
$ws.Range("G5").Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G6").Value = "Dr. Mohammad El-Tanany, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Alshimaa Atef"
$ws.Range("G7").Value = "Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G8").Value = "Dr. Nada Mohammad, Dr. Abeer Ragab"
$ws.Range("G11").Value = "Dr. Safa Hany, Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Range("G12").Value = "Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G13").Value = "Dr. Amira Ibrahim, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh"
$ws.Range("G15").Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef"
$ws.Range("G17").Value = "Dr. Mohammad Safwat, Dr. Esraa Samy"
$ws.Range("G19").Value = "Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef"
$ws.Range("G20").Value = "Dr. Mohammad Safwat, Dr. Mariam Toma Gerges"
$ws.Range("G25").Value = "Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil"
$ws.Range("G30").Value = "Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad"
